$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, shifting existing rows 145:203 down to 146:204.
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new weekly record.
$ws.Cells.Item(145, 1).Value = 6
$ws.Cells.Item(145, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(145, 3).Value = "Metropolitana"
$ws.Cells.Item(145, 4).Value = 44704
$ws.Cells.Item(145, 5).Value = 13
$ws.Cells.Item(145, 6).Value = 100112029
$ws.Cells.Item(145, 7).Value = "Orégano"
$ws.Cells.Item(145, 8).Value = "Sin especificar"
$ws.Cells.Item(145, 9).Value = "Primera"
$ws.Cells.Item(145, 10).Value = 34
$ws.Cells.Item(145, 11).Value = 12000
$ws.Cells.Item(145, 12).Value = 13000
$ws.Cells.Item(145, 13).Value = 12441
$ws.Cells.Item(145, 14).Value = "$/docena de atados"
$ws.Cells.Item(145, 15).Value = "Región Metropolitana"
$ws.Cells.Item(145, 16).Value = 4147
$ws.Cells.Item(145, 17).Value = 3
$ws.Cells.Item(145, 18).Value = "Hortaliza"
